# Update "想去人数" (want-to-go count) figures to the newly scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5252
$ws1.Range("F3").Value = 378
$ws1.Range("F6").Value = 796
$ws1.Range("F7").Value = 287
$ws1.Range("F8").Value = 10

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 7

# Sheet "全部类型" (All types, combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5252
$ws4.Range("F3").Value = 378
$ws4.Range("F6").Value = 796
$ws4.Range("F8").Value = 287
$ws4.Range("F9").Value = 10
$ws4.Range("F10").Value = 7
